# TC07 - Add to Wishlist error
# Adds a new "TC07" sheet (after the existing TC06 tab) that holds the
# baseUrl + the "Add to Wishlist" error-message fixture used by the test.

$wb = $excel.ActiveWorkbook

# Reset the previously-active sheet's selection back to A1 before we move
# away from it (mirrors the author switching tabs to create TC07).
$tc06 = $wb.Worksheets.Item("TC06")
$tc06.Activate()
$tc06.Range("A1").Select()

# Create the new sheet right after the last existing tab (TC06) and name it.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "TC07"

# Header row.
$ws.Range("A1").Value = "baseUrl"
$ws.Range("B1").Value = "errorMsg"

# Data row.
$ws.Range("A2").Value = "http://automationpractice.com/index.php"
$ws.Range("B2").Value = "You must be logged in to manage your wishlist."
$ws.Range("B2").WrapText = $true

# Column widths to fit the url / long error message.
$ws.Columns.Item(1).ColumnWidth = 33.94
$ws.Columns.Item(2).ColumnWidth = 40.01

# Match the page setup used by the rest of the workbook's sheets.
$ws.PageSetup.LeftMargin = 56.7
$ws.PageSetup.RightMargin = 56.7
$ws.PageSetup.TopMargin = 75.8
$ws.PageSetup.BottomMargin = 75.8
$ws.PageSetup.HeaderMargin = 56.7
$ws.PageSetup.FooterMargin = 56.7
$ws.PageSetup.CenterHeader = '&"Times New Roman,Normal"&12&A'
$ws.PageSetup.CenterFooter = '&"Times New Roman,Normal"&12Página &P'

# Leave the selection where the author left it and make sure TC07 is active.
$ws.Range("D8").Select()
$ws.Activate()
